$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last refreshed" timestamp banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Mayo de 2020 a las 12:08"

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Countries that swapped ranking order (Bielorrusia overtook Austria; ---
# --- Indonesia overtook Corea del Sur) together with refreshed case data ---

# Row 31 now holds Bielorrusia (its case counts grew past Austria's)
Set-Row 31 @("Bielorrusia", 15828, 911, 3117, 12614, 92, 4, 97)

# Row 32 now holds Austria, keeping its previous (unrefreshed) case data
Set-Row 32 @("Austria", 15531, 0, 13110, 1832, 124, 0, 589)

# Row 39 now holds Indonesia (its case counts grew past Corea del Sur's)
Set-Row 39 @("Indonesia", 10843, 292, 1665, 8347, 0, 31, 831)

# Row 40 now holds Corea del Sur, keeping its previous (unrefreshed) case data
Set-Row 40 @("Corea del Sur", 10780, 6, 9123, 1407, 55, 2, 250)

# --- Plain data refreshes (country stays in the same row) ---

# Belgica
Set-Row 16 @("Belgica", 49517, 485, 12211, 29541, 689, 62, 7765)

# Suiza
Set-Row 20 @("Suiza", 29817, 112, 23900, 4163, 167, 0, 1754)

# Rumania
Set-Row 37 @("Rumania", 12567, 0, 4328, 7484, 249, 11, 755)

# Finlandia
Set-Row 54 @("Finlandia", 5176, 125, 3000, 1958, 49, 0, 218)

# Kazajistan
Set-Row 61 @("Kazajistan", 3785, 188, 940, 2820, 40, 0, 25)

# Afganistan
Set-Row 67 @("Afganistan", 2469, 134, 331, 2066, 7, 4, 72)

# Albania
Set-Row 96 @("Albania", 789, 7, 519, 239, 4, 0, 31)

# Etiopia
Set-Row 143 @("Etiopia", 133, 0, 69, 61, 0, 0, 3)
